$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J (10) - shifts old J/K to K/L.
$ws.Columns.Item(10).Insert()

# New header cell J1 "Korrektur" (re-uses header style from the row, same
# vertical-top style as the other un-wrapped headers).
$ws.Range("J1").Value = "Korrektur"

# New data cell J2 "{isKorrektur}" with currency numFmt + right alignment,
# matching the new cellXf that Excel creates for this cell.
$ws.Range("J2").Value = "{isKorrektur}"
$ws.Range("J2").NumberFormat = '[$CHF]\ #,##0.00;[Red]\-[$CHF]\ #,##0.00'
$ws.Range("J2").HorizontalAlignment = -4152

# Match column J's width to column I's width (both become the
# "13.42578125" width class used throughout the template).
$ws.Columns.Item(10).ColumnWidth = $ws.Columns.Item(9).ColumnWidth

# Highlight whole rows in red font when the row is flagged as a
# correction (J column holds "X").
$rng = $ws.Range("A2:XFD99999")
$fc = $rng.FormatConditions.Add(2, 0, '$J$2="X"')
$fc.Font.Color = 255

$wb.Save()
